$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Row 3: Wald (C,E,G,I) and bitmap (D,F,H,J)
$ws.Range("C3").Value = "Wald, Gras"
$ws.Range("E3").Value = "Wald, Gras"
$ws.Range("G3").Value = "Wald, Gras"
$ws.Range("I3").Value = "Wald, Gras"
$ws.Range("D3").Value = "000000011"
$ws.Range("F3").Value = "000000011"
$ws.Range("H3").Value = "000000011"
$ws.Range("J3").Value = "000000011"

# Row 4: Kuh (C,E,G,I) and bitmap (D,F,H,J)
$ws.Range("C4").Value = "Gras"
$ws.Range("E4").Value = "Gras"
$ws.Range("G4").Value = "Gras"
$ws.Range("I4").Value = "Gras"
$ws.Range("D4").Value = "000000001"
$ws.Range("F4").Value = "000000001"
$ws.Range("H4").Value = "000000001"
$ws.Range("J4").Value = "000000001"

# Row 8: Berg (C,E,G,I) and bitmap (D,F,H,J)
$ws.Range("C8").Value = "Gras, Berg, Schnee"
$ws.Range("E8").Value = "Gras, Berg, Schnee"
$ws.Range("G8").Value = "Gras, Berg, Schnee"
$ws.Range("I8").Value = "Gras, Berg, Schnee"
$ws.Range("D8").Value = "011000001"
$ws.Range("F8").Value = "011000001"
$ws.Range("H8").Value = "011000001"
$ws.Range("J8").Value = "011000001"

# Row 5: Strand (C,E,G,I) and bitmap (D,F,H,J)
$ws.Range("C5").Value = "Wasser, Gras, Strand"
$ws.Range("E5").Value = "Wasser, Gras, Strand"
$ws.Range("G5").Value = "Wasser, Gras, Strand"
$ws.Range("I5").Value = "Wasser, Gras, Strand"
$ws.Range("D5").Value = "000011001"
$ws.Range("F5").Value = "000011001"
$ws.Range("H5").Value = "000011001"
$ws.Range("J5").Value = "000011001"

# Update selection to reflect the active cell used while editing
$ws.Range("C3").Select()

$wb.Save()
